$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 218.42857
$ws.Range("I33").Value = 130.66667
$ws.Range("J33").Value = 745
$ws.Range("K33").Value = 130.66667
$ws.Range("L33").Value = 745
$ws.Range("M33").Value = 98.33332999999999
$ws.Range("N33").Value = -1203

$ws.Range("H64").Value = 3364.0667
$ws.Range("I64").Value = 3342.923
$ws.Range("K64").Value = 3342.923
$ws.Range("M64").Value = -3094.923

$ws.Range("H67").Value = 3364.0667
$ws.Range("I67").Value = 3342.923
$ws.Range("K67").Value = 3342.923
$ws.Range("M67").Value = -2484.923

$ws.Range("H74").Value = 4948.2666
$ws.Range("I74").Value = 4948.2666
$ws.Range("K74").Value = 4948.2666
$ws.Range("M74").Value = -4012.2666

$ws.Range("H77").Value = 4948.2666
$ws.Range("I77").Value = 4948.2666
$ws.Range("K77").Value = 24741.333
$ws.Range("M77").Value = -20061.333

$ws.Range("H99").Value = 115079660
$ws.Range("I99").Value = 5102378.5
$ws.Range("K99").Value = 15307135.5
$ws.Range("M99").Value = -15305637.5

$ws.Range("H100").Value = 4897.273
$ws.Range("I100").Value = 2474.5557
$ws.Range("K100").Value = 2474.5557
$ws.Range("M100").Value = -1933.5557

$ws.Range("H115").Value = 6175430.5
$ws.Range("I115").Value = 6175430.5
$ws.Range("K115").Value = 18526291.5
$ws.Range("M115").Value = -18524724.5

$ws.Range("H116").Value = 4300.7666
$ws.Range("I116").Value = 3979.3333
$ws.Range("J116").Value = 4782.9165
$ws.Range("K116").Value = 3979.3333
$ws.Range("L116").Value = 4782.9165
$ws.Range("M116").Value = -537.3332999999998
$ws.Range("N116").Value = -11666.9165

$ws.Range("H118").Value = 3249123
$ws.Range("I118").Value = 4465268
$ws.Range("K118").Value = 13395804
$ws.Range("M118").Value = -13394147

$ws.Range("H129").Value = 1164.25
$ws.Range("I129").Value = 1164.25
$ws.Range("K129").Value = 3492.75
$ws.Range("M129").Value = 1507.25

$ws.Range("H137").Value = 52142.184
$ws.Range("I137").Value = 79863.42999999999
$ws.Range("K137").Value = 239590.29
$ws.Range("M137").Value = -237040.29

$ws.Range("H138").Value = 3202.8718
$ws.Range("J138").Value = 3541.4333
$ws.Range("L138").Value = 10624.2999
$ws.Range("N138").Value = -20904.2999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3776.2222
$ws.Range("I63").Value = 3776.2222
$ws.Range("K63").Value = 3776.2222
$ws.Range("M63").Value = -3090.2222

$ws.Range("H66").Value = 3776.2222
$ws.Range("I66").Value = 3776.2222
$ws.Range("K66").Value = 18881.111
$ws.Range("M66").Value = -15449.111

$ws.Range("H74").Value = 493191
$ws.Range("I74").Value = 656921.7
$ws.Range("J74").Value = 1999
$ws.Range("K74").Value = 656921.7
$ws.Range("L74").Value = 1999
$ws.Range("M74").Value = -656047.7
$ws.Range("N74").Value = -3747

$ws.Range("H77").Value = 493191
$ws.Range("I77").Value = 656921.7
$ws.Range("J77").Value = 1999
$ws.Range("K77").Value = 3284608.5
$ws.Range("L77").Value = 9995
$ws.Range("M77").Value = -3280240.5
$ws.Range("N77").Value = -18731

$ws.Range("H97").Value = 898.38776
$ws.Range("I97").Value = 655.9459000000001
$ws.Range("J97").Value = 1645.9166
$ws.Range("K97").Value = 655.9459000000001
$ws.Range("L97").Value = 1645.9166
$ws.Range("M97").Value = -159.9459000000001
$ws.Range("N97").Value = -2637.9166

$ws.Range("H122").Value = 3240.5715
$ws.Range("I122").Value = 1537.5
$ws.Range("K122").Value = 4612.5
$ws.Range("M122").Value = -2162.5

$ws.Range("H132").Value = 3835.2222
$ws.Range("I132").Value = 3871.5
$ws.Range("K132").Value = 11614.5
$ws.Range("M132").Value = -9084.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 708.5454999999999
$ws.Range("I22").Value = 708.5454999999999
$ws.Range("K22").Value = 708.5454999999999
$ws.Range("M22").Value = -535.5454999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("I86").Value = 1433842.8
$ws.Range("J86").Value = 4374
$ws.Range("K86").Value = 1433842.8
$ws.Range("L86").Value = 4374
$ws.Range("M86").Value = -1432719.8
$ws.Range("N86").Value = -6620

$ws.Range("I89").Value = 1433842.8
$ws.Range("J89").Value = 4374
$ws.Range("K89").Value = 7169214
$ws.Range("L89").Value = 21870
$ws.Range("M89").Value = -7163598
$ws.Range("N89").Value = -33102

$ws.Range("H132").Value = 7022.9165
$ws.Range("I132").Value = 4475.4443
$ws.Range("K132").Value = 13426.3329
$ws.Range("M132").Value = -10896.3329

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 8246.388999999999
$ws.Range("I132").Value = 10954.333
$ws.Range("K132").Value = 98588.997
$ws.Range("M132").Value = -96058.997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 502512
$ws.Range("I132").Value = 502512
$ws.Range("K132").Value = 1507536
$ws.Range("M132").Value = -1505006

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 930.6667
$ws.Range("I22").Value = 1000
$ws.Range("J22").Value = 896
$ws.Range("K22").Value = 1000
$ws.Range("L22").Value = 896
$ws.Range("M22").Value = -705
$ws.Range("N22").Value = -1486

$ws.Range("H27").Value = 930.6667
$ws.Range("I27").Value = 1000
$ws.Range("J27").Value = 896
$ws.Range("K27").Value = 1000
$ws.Range("L27").Value = 1000
$ws.Range("M27").Value = -893
$ws.Range("N27").Value = -1110

$ws.Range("H43").Value = 19111.611
$ws.Range("J43").Value = 19145.727
$ws.Range("N43").Value = -19531.727

$ws.Range("H122").Value = 10000
$ws.Range("I122").Value = 10000
$ws.Range("K122").Value = 30000
$ws.Range("M122").Value = -27550

$ws.Range("H132").Value = 4301.7915
$ws.Range("I132").Value = 4075.0557
$ws.Range("K132").Value = 12225.1671
$ws.Range("M132").Value = -9695.167099999999

$ws.Range("H136").Value = 2441.2
$ws.Range("I136").Value = 1708.85
$ws.Range("K136").Value = 5126.549999999999
$ws.Range("M136").Value = -2576.549999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H59").Value = 24999
$ws.Range("J59").Value = 24999
$ws.Range("L59").Value = 24999
$ws.Range("N59").Value = -26475

$ws.Range("H96").Value = 5138.909
$ws.Range("I96").Value = 5005.1665
$ws.Range("K96").Value = 5005.1665
$ws.Range("M96").Value = -3632.1665

$ws.Range("H122").Value = 4051
$ws.Range("I122").Value = 4051
$ws.Range("K122").Value = 12153
$ws.Range("M122").Value = -9703
